# Applies the "DB finish constrains, user cases and zugriffsmatrix" edit:
#  - anr/wnr/knr/bnr/annr/bildnr "Datentyp" columns switch from String(3) to Integer
#    (the Laenge/E column is cleared since Integer has no declared length here)
#  - several foreign-key Constraint cells switch ON DELETE RESTRICT -> ON DELETE CASCADE
#  - Ausstattung.name (H19) becomes a primary-key flag ("j")
#  - iban description loses its stray leading space
#  - Bewertung/Rechnung rows: Muss/Kann flips from "m" to "k", constraint text tightened,
#    and a new ">= 0" constraint is added for Rechnungsbetrag
#  - Anzahlung.bnr / Anzahlung.datum drop their "Primaerschluessel" flag (j -> n)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Datentyp: String -> Integer, and clear the now-unused Laenge (E) value ---
$intRows = @(5, 16, 21, 24, 25, 27, 28, 33, 34, 40, 41, 42, 49, 52, 53)
foreach ($r in $intRows) {
    $ws.Range("D$r").Value = "Integer"
    $ws.Range("E$r").ClearContents()
}

# --- Constraint: ON DELETE RESTRICT -> ON DELETE CASCADE ---
$cascadeRows = @(16, 17, 21, 22, 25)
foreach ($r in $cascadeRows) {
    $ws.Range("J$r").Value = "ON DELETE CASCADE"
}

# --- Ausstattung.name becomes a primary key ---
$ws.Range("H19").Value = "j"

# --- iban description: drop the leading space ---
$ws.Range("C38").Value = "Bankverbindung"

# --- Bewertung constraint text tightened ---
$ws.Range("J48").Value = "1 <= anz <= 5"

# --- Rechnung rows: Muss/Kann m -> k ---
$ws.Range("F49").Value = "k"
$ws.Range("F50").Value = "k"
$ws.Range("F51").Value = "k"

# --- new constraint on Rechnungsbetrag ---
$ws.Range("J51").Value = ">= 0"

# --- Anzahlung.bnr / Anzahlung.datum no longer flagged as primary key ---
$ws.Range("H53").Value = "n"
$ws.Range("H54").Value = "n"

# --- cosmetic: new helper columns (F / K) sized like the other narrow columns ---
$ws.Columns.Item(6).ColumnWidth = 10
$ws.Columns.Item(11).ColumnWidth = 10

# --- page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- selection left where the author ended up ---
$ws.Range("N54").Select() | Out-Null
